$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header text for the MODEL_CONDITION column (now "MODELCONDITION")
$ws.Range("E1").Value = "MODELCONDITION"

# Delete entire column A (the old TAXON-number-like column), shifting B:F left to A:E
$ws.Columns("A").Delete()
